$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the bootstrap note to the existing journal entry in E34
$ws.Range("E34").Value = "Still running grid search, added back in grades < 60. Removing rows without prereqs helped a bunch on train/test. Quantile error bars look good w new data. Bootstrap error bars on xgboost"

# Hours worked that day went from 4 to 8
$ws.Range("C34").Value = 8

# Row 34 grows from 3 wrapped lines to 4, so the row gets taller
$ws.Rows.Item(34).RowHeight = 57

# New TODO note added in row 35
$ws.Range("G35").Value = "should play around with bootstrapping more"

# Move the selection/cursor to reflect where the author ended up after typing
$ws.Range("G36").Select() | Out-Null
